# Update forecast values on the "Forecast Comparison" sheet (Removed Auto Arima)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

$data = @(
    @(13, 17, 20, 26),
    @(13, 16, 21, 27),
    @(14, 18, 22, 29),
    @(15, 18, 23, 31),
    @(15, 19, 24, 32),
    @(15, 19, 23, 31),
    @(16, 19, 25, 33),
    @(16, 20, 25, 34),
    @(15, 19, 24, 32),
    @(15, 19, 24, 33),
    @(15, 19, 25, 35),
    @(17, 20, 26, 36),
    @(15, 18, 24, 34),
    @(14, 18, 24, 33),
    @(14, 17, 23, 32),
    @(13, 15, 21, 30)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 4).Value = $vals[0]
    $ws.Cells.Item($row, 5).Value = $vals[1]
    $ws.Cells.Item($row, 6).Value = $vals[2]
    $ws.Cells.Item($row, 7).Value = $vals[3]
}
